# Scheduled runner update: refresh market-price derived columns (H-N) per row
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 250
$ws.Range("I5").Value = 278.66666
$ws.Range("J5").Value = 78
$ws.Range("K5").Value = 278.66666
$ws.Range("L5").Value = 78
$ws.Range("M5").Value = -163.66666
$ws.Range("N5").Value = -308

$ws.Range("H9").Value = 119.6
$ws.Range("I9").Value = 50.142857
$ws.Range("J9").Value = 281.66666
$ws.Range("K9").Value = 50.142857
$ws.Range("L9").Value = 281.66666
$ws.Range("M9").Value = 118.857143
$ws.Range("N9").Value = -619.66666

$ws.Range("H12").Value = 343.4
$ws.Range("I12").Value = 378.75
$ws.Range("J12").Value = 202
$ws.Range("K12").Value = 378.75
$ws.Range("L12").Value = 202
$ws.Range("M12").Value = -208.75
$ws.Range("N12").Value = -542

$ws.Range("H53").Value = 198
$ws.Range("I53").Value = 163.33333
$ws.Range("J53").Value = 218.8
$ws.Range("K53").Value = 163.33333
$ws.Range("L53").Value = 218.8
$ws.Range("M53").Value = 473.66667
$ws.Range("N53").Value = -1492.8

$ws.Range("H76").Value = 3000
$ws.Range("J76").Value = 3000
$ws.Range("L76").Value = 3000
$ws.Range("N76").Value = -3630

$ws.Range("H79").Value = 3000
$ws.Range("J79").Value = 3000
$ws.Range("L79").Value = 3000
$ws.Range("N79").Value = -5184

$ws.Range("H132").Value = 2449
$ws.Range("I132").Value = 786.25
$ws.Range("K132").Value = 2358.75
$ws.Range("M132").Value = 171.25

$ws.Range("H138").Value = 2999
$ws.Range("I138").Value = 2999
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 8997
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -3857
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 107
$ws.Range("I5").Value = 103
$ws.Range("K5").Value = 103
$ws.Range("M5").Value = 9

$ws.Range("H61").Value = 500
$ws.Range("I61").Value = 500
$ws.Range("K61").Value = 500
$ws.Range("M61").Value = -288

$ws.Range("H132").Value = 883.3333
$ws.Range("I132").Value = 883.3333
$ws.Range("K132").Value = 2649.9999
$ws.Range("M132").Value = -119.9998999999998

$ws.Range("H136").Value = 500
$ws.Range("I136").Value = 500
$ws.Range("K136").Value = 1500
$ws.Range("M136").Value = 1050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 107
$ws.Range("I4").Value = 103
$ws.Range("K4").Value = 103
$ws.Range("M4").Value = 12

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H134").Value = 9244.666999999999
$ws.Range("I134").Value = 7628.4614
$ws.Range("K134").Value = 22885.3842
$ws.Range("M134").Value = -20350.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 556.8182
$ws.Range("J22").Value = 400
$ws.Range("L22").Value = 400
$ws.Range("N22").Value = -1100

$ws.Range("H99").Value = 852437.5
$ws.Range("I99").Value = 579944.9
$ws.Range("J99").Value = 1669915.4
$ws.Range("K99").Value = 579944.9
$ws.Range("L99").Value = 1669915.4
$ws.Range("M99").Value = -578446.9
$ws.Range("N99").Value = -1672911.4

$ws.Range("H122").Value = 1126.75
$ws.Range("I122").Value = 836
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 2508
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -58
$ws.Range("N122").Value = -10897

$ws.Range("H126").Value = 852437.5
$ws.Range("I126").Value = 579944.9
$ws.Range("J126").Value = 1669915.4
$ws.Range("K126").Value = 1739834.7
$ws.Range("L126").Value = 5009746.199999999
$ws.Range("M126").Value = -1737364.7
$ws.Range("N126").Value = -5014686.199999999

$ws.Range("H132").Value = 3403.5
$ws.Range("I132").Value = 3500
$ws.Range("K132").Value = 10500
$ws.Range("M132").Value = -7970

$ws.Range("H134").Value = 2587.3333
$ws.Range("I134").Value = 2587.3333
$ws.Range("K134").Value = 7761.999899999999
$ws.Range("M134").Value = -5226.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2648.2856
$ws.Range("I3").Value = 999
$ws.Range("K3").Value = 2997
$ws.Range("M3").Value = -2885

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2144.25
$ws.Range("I22").Value = 924.75
$ws.Range("J22").Value = 3363.75
$ws.Range("K22").Value = 924.75
$ws.Range("L22").Value = 3363.75
$ws.Range("M22").Value = -629.75
$ws.Range("N22").Value = -3953.75

$ws.Range("H27").Value = 2144.25
$ws.Range("I27").Value = 924.75
$ws.Range("J27").Value = 3363.75
$ws.Range("K27").Value = 924.75
$ws.Range("L27").Value = 3363.75
$ws.Range("M27").Value = -817.75
$ws.Range("N27").Value = -3577.75

$ws.Range("H40").Value = 17115.834
$ws.Range("I40").Value = 2538
$ws.Range("K40").Value = 2538
$ws.Range("M40").Value = -2402

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H122").Value = 7997.5
$ws.Range("I122").Value = 7995
$ws.Range("K122").Value = 23985
$ws.Range("M122").Value = -21535

$ws.Range("H136").Value = 2561.625
$ws.Range("I136").Value = 2658.6
$ws.Range("K136").Value = 7975.799999999999
$ws.Range("M136").Value = -5425.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 10003.333
$ws.Range("I30").Value = 10000
$ws.Range("J30").Value = 10005
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 10005
$ws.Range("M30").Value = -9893
$ws.Range("N30").Value = -10219

$ws.Range("H122").Value = 2660.7778
$ws.Range("I122").Value = 2387.5
$ws.Range("J122").Value = 2879.4
$ws.Range("K122").Value = 7162.5
$ws.Range("L122").Value = 8638.200000000001
$ws.Range("M122").Value = -4712.5
$ws.Range("N122").Value = -13538.2

$ws.Range("H132").Value = 2599.5
$ws.Range("I132").Value = 2599.5
$ws.Range("K132").Value = 7798.5
$ws.Range("M132").Value = -5268.5

$ws.Range("H136").Value = 4097.6924
$ws.Range("I136").Value = 522.8333
$ws.Range("K136").Value = 1568.4999
$ws.Range("M136").Value = 981.5001
